$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F (想去人数) values
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1058
$ws.Range("F5").Value = 19557
$ws.Range("F6").Value = 74
$ws.Range("F7").Value = 2207
$ws.Range("F10").Value = 430
$ws.Range("F13").Value = 249
$ws.Range("F19").Value = 186
$ws.Range("F22").Value = 97

# Sheet "演出" (Performances) - update column F (想去人数) values
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 192
$ws.Range("F4").Value = 10
$ws.Range("F5").Value = 15
$ws.Range("F14").Value = 38
$ws.Range("F21").Value = 34

# Sheet "全部类型" (All types) - update column F (想去人数) values
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 192
$ws.Range("F9").Value = 1058
$ws.Range("F10").Value = 19557
$ws.Range("F11").Value = 10
$ws.Range("F12").Value = 15
$ws.Range("F13").Value = 74
$ws.Range("F16").Value = 2207
$ws.Range("F20").Value = 430
$ws.Range("F23").Value = 249
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 186
$ws.Range("F44").Value = 34
$ws.Range("F47").Value = 97
